# Applies the USCDI table edit:
#  1. Insert a new column at H ("fhir_extension_url"), which pushes the
#     existing fhir_path/comment/comment_link (and every other trailing
#     column) one slot to the right. This also naturally extends the
#     used range so a new (empty) column shows up at the end (S).
#  2. Set the new header text in H1.
#  3. Populate the few rows that now carry a real
#     fhir_extension_url value (StructureDefinition-*.html links).
#  4. Apply the handful of text-content tweaks to uscore_profile /
#     extension columns called out by the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert new column before column H (8) ---------------------------
$ws.Columns.Item(8).Insert()

# --- 2. New header label for the inserted column -------------------------
$ws.Range("H1").Value = "fhir_extension_url"

# --- 3. New fhir_extension_url values (StructureDefinition links) -------
$ws.Range("H90").Value  = "StructureDefinition-patient-sexParameterForClinicalUse.html"
$ws.Range("H104").Value = "StructureDefinition-individual-pronouns.html"
$ws.Range("H127").Value = "StructureDefinition-condition-assertedDate.html"

# --- 4. uscore_profile (C) / extension (G) / fhir_path (I) text tweaks ---
$ws.Range("C86").Value  = "US Core MedicationRequest Profile | US Core Medication Adherence Extension"

$ws.Range("C90").Value  = "US Core Patient Profile | Patient Sex Parameter for Clinical Use Extension"

$ws.Range("C104").Value = "US Core Patient Profile | Individual Pronouns Extension"

$ws.Range("C108").Value = "US Core Patient Profile | US Core Race Extension"

$ws.Range("C114").Value = "US Core Patient Profile | US Core Encounter Profile |US Core Interpreter Needed Extension"
$ws.Range("G114").Value = "US Core Interpreter Needed Extension"
$ws.Range("I114").Value = "Patient.extension.where(url='http://hl7.org/fhir/us/core/StructureDefinition/us-core-interpreter-needed'') | Encounter.extension.where(url='http://hl7.org/fhir/us/core/StructureDefinition/us-core-interpreter-needed'') "

$ws.Range("C126").Value = "US Core Condition Problems and Health Concerns Profile |assertedDate Extension"

$ws.Range("C127").Value = "US Core Condition Problems and Health Concerns Profile |assertedDate Extension"
